$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Paris 2023 Legends Sticker Capsule"
$ws.Range("A3").Value = "Paris 2023 Challengers Sticker Capsule"
$ws.Range("A4").Value = "Dreams & Nightmares Case"
